$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.893.67"
Set-TextValue "E2" "  -2.26%  "

Set-TextValue "D3" "1.779.00"
Set-TextValue "E3" "  -2.86%  "

Set-TextValue "D4" "1.010"
Set-TextValue "E4" "  +0.82%  "

Set-TextValue "D5" "1.010"
Set-TextValue "E5" "  +0.87%  "

Set-TextValue "D6" "308.71"
Set-TextValue "E6" "  -1.20%  "

Set-TextValue "D7" "0.4230"
Set-TextValue "E7" "  -1.60%  "

Set-TextValue "D8" "0.3618"
Set-TextValue "E8" "  -1.15%  "

Set-TextValue "D9" "0.07178"
Set-TextValue "E9" "  -1.20%  "

Set-TextValue "D10" "0.8371"
Set-TextValue "E10" "  -3.53%  "

Set-TextValue "D11" "20.24"
Set-TextValue "E11" "  -1.97%  "

Set-TextValue "D12" "1.826.82"
Set-TextValue "E12" "  +0.11%  "

Set-TextValue "D13" "5.244"
Set-TextValue "E13" "  -2.84%  "

Set-TextValue "D14" "6.332"
Set-TextValue "E14" "  -2.97%  "

Set-TextValue "D15" "0.06819"
Set-TextValue "E15" "  -1.68%  "

Set-TextValue "D16" "1.014"
Set-TextValue "E16" "  +1.23%  "

Set-TextValue "D17" "79.09"
Set-TextValue "E17" "  -1.73%  "

Set-TextValue "D18" "0.000008666"
Set-TextValue "E18" "  -2.74%  "

Set-TextValue "D19" "1.010"
Set-TextValue "E19" "  +0.84%  "

Set-TextValue "D20" "14.93"
Set-TextValue "E20" "  -2.88%  "

Set-TextValue "D21" "27.158.40"
Set-TextValue "E21" "  -1.17%  "

Set-TextValue "D22" "5.015"
Set-TextValue "E22" "  -2.26%  "

Set-TextValue "D23" "11.03"
Set-TextValue "E23" "  +1.86%  "

Set-TextValue "D24" "2.054.69"
Set-TextValue "E24" "  +0.35%  "

Set-TextValue "D25" "1.926"
Set-TextValue "E25" "  -2.64%  "

Set-TextValue "D26" "153.49"
Set-TextValue "E26" "  -0.57%  "

Set-TextValue "D27" "18.12"
Set-TextValue "E27" "  -4.10%  "

Set-TextValue "B28" "InternetComputer(DFINITY)"
Set-TextValue "C28" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D28" "5.025"
Set-TextValue "E28" "  -2.38%  "

Set-TextValue "B29" "BitcoinCash"
Set-TextValue "C29" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D29" "114.07"
Set-TextValue "E29" "  -0.04%  "

Set-TextValue "D30" "1.628"
Set-TextValue "E30" "  -11.13%  "

Set-TextValue "D31" "0.08943"
Set-TextValue "E31" "  +0.78%  "

Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "0.7154"
Set-TextValue "E32" "  -5.03%  "

Set-TextValue "B33" "HuobiToken"
Set-TextValue "C33" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D33" "2.855"
Set-TextValue "E33" "  -4.26%  "

Set-TextValue "D34" "4.315"
Set-TextValue "E34" "  -4.99%  "

Set-TextValue "D35" "1.086"
Set-TextValue "E35" "  -4.20%  "

Set-TextValue "D36" "1.011"
Set-TextValue "E36" "  +0.93%  "

Set-TextValue "D37" "1.079"
Set-TextValue "E37" "  -1.02%  "

Set-TextValue "D38" "0.01905"
Set-TextValue "E38" "  -1.52%  "

Set-TextValue "D39" "0.05074"
Set-TextValue "E39" "  -4.50%  "

Set-TextValue "B40" "Algorand"
Set-TextValue "C40" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D40" "0.1607"
Set-TextValue "E40" "  -3.69%  "

Set-TextValue "B41" "TheSandbox"
Set-TextValue "C41" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.4905"
Set-TextValue "E41" "  -3.41%  "

Set-TextValue "D42" "2.499"
Set-TextValue "E42" "  -10.74%  "

Set-TextValue "D43" "6.022"
Set-TextValue "E43" "  -8.57%  "

Set-TextValue "D44" "7.902"
Set-TextValue "E44" "  -5.79%  "

Set-TextValue "D45" "1.010"
Set-TextValue "E45" "  +1.03%  "

Set-TextValue "D46" "104.31"
Set-TextValue "E46" "  -1.64%  "

Set-TextValue "D47" "10.03"
Set-TextValue "E47" "  -3.85%  "

Set-TextValue "D48" "0.06227"
Set-TextValue "E48" "  -4.14%  "

Set-TextValue "D49" "0.4446"
Set-TextValue "E49" "  -4.96%  "

Set-TextValue "D50" "1.568"
Set-TextValue "E50" "  -2.81%  "

Set-TextValue "D51" "1.688"
Set-TextValue "E51" "  -2.79%  "
